# Apply the OOXML diff described:
#  - Sheet "Metadata": Date value B8 changes from 2025-07-18T06:40:38+00:00 to 2025-07-21T11:52:46+00:00
#  - Sheet "Include #0": System URI B4 changes from ansforge.github.io ... R01 ... to interop.esante.gouv.fr ... R01 ...
#  - Sheet "Include #1": System URI B4 changes from ansforge.github.io ... R359 ... to interop.esante.gouv.fr ... R359 ...

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B8").Value = "2025-07-21T11:52:46+00:00"

$wsInclude0 = $wb.Worksheets.Item("Include #0")
$wsInclude0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R01-EnsembleSavoirFaire-CISIS"

$wsInclude1 = $wb.Worksheets.Item("Include #1")
$wsInclude1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R359-SurspecialiteTransversale"
